$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 swap their Name (B) / Class (D) values.
$ws.Range("B2").Value = "Tô Hiếu Ngôi"
$ws.Range("D2").Value = "59.CNTT-3"

$ws.Range("B3").Value = "Lê Thế Dũng"
$ws.Range("D3").Value = "59.CNTT-1"

# Row 6: replace the student's name and update their class.
$ws.Range("B6").Value = "Nguyễn Xuân Huy"
$ws.Range("D6").Value = "59.CNTT-3"

# Move the active selection to D13 (was B13).
$ws.Range("D13").Select()
